# GAP-10 added geoserver data to config file
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the existing "dbgap" parameter value to "db_gap"
$ws.Range("B5").Value = "db_gap"

# Append new geoserver-related parameter rows below the existing table
$ws.Range("A7").Value = "user_gs"

$ws.Range("A8").Value = "pass_gs"
$ws.Range("B8").Style = "Normal"

$ws.Range("A9").Value = "workspace_gs"

$ws.Range("A10").Value = "path"

$ws.Range("A11").Value = "geo_url"

# B11 (the geo_url value cell) gets hyperlink-style formatting, left blank for now
$b11 = $ws.Range("B11")
$b11.Style = "Normal"
$ws.Hyperlinks.Add($b11, "https://example.com")
$ws.Hyperlinks.Item(1).Delete()
$b11.ClearContents()

# Restore selection and page setup as saved
[void]$ws.Range("B5").Select()
$ws.PageSetup.Orientation = 1
